$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order: old row 4, old row 5, old row 3, old row 2
# Only columns D, L, M, N, O, P, S change; everything else stays identical.

$ws.Range("D2").Value = 45008
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("S2").Value = 3500

$ws.Range("D3").Value = 45008
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("S3").Value = 3000

$ws.Range("D4").Value = 44991
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000

$ws.Range("D5").Value = 44995
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5500
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 5750
$ws.Range("S5").Value = 2875
